## Update Data by bot, scripted by HH
# Row 2 of Sheet1 is a single financial-data record (300925.SZ). This
# script refreshes it from the 2020-06-30 report to the 2017-12-31 report,
# updating the date-type code, report date and the associated balance-
# sheet figures (and filling in two previously-blank ratio/amount cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE ("002" -> "001"). Force-write as text (apostrophe prefix
# keeps the leading zero instead of Excel coercing it to the number 1),
# then reset the cell style back to Normal so no stray number-format /
# quote-prefix styling is left behind on the cell.
$ws.Range("J2").Value = "'001"
$ws.Range("J2").Style = "Normal"

# REPORT_DATE
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# TOTAL_ASSETS / FIXED_ASSET / MONETARYFUNDS / MONETARYFUNDS_RATIO
$ws.Range("O2").Value = 255745745.84
$ws.Range("P2").Value = 3517204.87
$ws.Range("Q2").Value = 36864383.82
$ws.Range("R2").Value = -53.4232197125

# ACCOUNTS_RECE / ACCOUNTS_RECE_RATIO
$ws.Range("S2").Value = 140674934.29
$ws.Range("T2").Value = 48.7144504492

# TOTAL_LIABILITIES / ACCOUNTS_PAYABLE / ACCOUNTS_PAYABLE_RATIO
$ws.Range("W2").Value = 84420860.90000001
$ws.Range("X2").Value = 3773010.74
$ws.Range("Y2").Value = 132.9863637168

# ADVANCE_RECEIVABLES / ADVANCE_RECEIVABLES_RATIO (previously blank)
$ws.Range("Z2").Value = 289760.28
$ws.Range("AA2").Value = -64.6233008448

# TOTAL_EQUITY / TOTAL_EQUITY_RATIO
$ws.Range("AB2").Value = 171324884.94
$ws.Range("AC2").Value = 73.63705071859999

# TOTAL_ASSETS_RATIO / TOTAL_LIAB_RATIO / CURRENT_RATIO / DEBT_ASSET_RATIO
$ws.Range("AD2").Value = 38.4444189546
$ws.Range("AE2").Value = -1.9043155876
$ws.Range("AF2").Value = 295.5228509995
$ws.Range("AG2").Value = 33.0096833567
